$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 19 to accommodate 2 net-new profile entries
# (list grows from 55 to 57 total rows: "US Core Implantable Device Profile" is
# removed while "US Core Device Profile", "US Core FamilyMemberHistory Profile" and
# "US Core PMO ServiceRequest Profile" are newly inserted in alphabetical order)
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Insert()

# Data for rows 19-57: row, index(col A), profile name(col B), target1..target6 (col C-H)
$rowsData = @(
    ,@(19, 17, "US Core Device Profile", $null, $null, $null, $null, $null, $null)
    ,@(20, 18, "US Core DiagnosticReport Profile for Laboratory Results Reporting", $null, $null, $null, $null, $null, $null)
    ,@(21, 19, "US Core DiagnosticReport Profile for Report and Note Exchange", $null, $null, $null, $null, $null, $null)
    ,@(22, 20, "US Core DocumentReference Profile", $null, $null, $null, $null, $null, $null)
    ,@(23, 21, "US Core Encounter Profile", $null, $null, $null, $null, $null, $null)
    ,@(24, 22, "US Core FamilyMemberHistory Profile", $null, $null, $null, $null, $null, $null)
    ,@(25, 23, "US Core Goal Profile", $null, $null, $null, $null, $null, $null)
    ,@(26, 24, "US Core Head Circumference Profile", "PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")
    ,@(27, 25, "US Core Heart Rate Profile", "PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")
    ,@(28, 26, "US Core Immunization Profile", $null, $null, $null, $null, $null, $null)
    ,@(29, 27, "US Core Location Profile", $null, $null, $null, $null, $null, $null)
    ,@(30, 28, "US Core Medication Profile", $null, $null, $null, $null, $null, $null)
    ,@(31, 29, "US Core MedicationDispense Profile", $null, $null, $null, $null, $null, $null)
    ,@(32, 30, "US Core MedicationRequest Profile", $null, $null, $null, $null, $null, $null)
    ,@(33, 31, "US Core Observation ADI Documentation Profile", $null, $null, $null, $null, $null, $null)
    ,@(34, 32, "US Core Observation Clinical Result Profile", $null, $null, $null, $null, $null, $null)
    ,@(35, 33, "US Core Laboratory Result Observation Profile", $null, $null, $null, $null, $null, $null)
    ,@(36, 34, "US Core Observation Occupation Profile", $null, $null, $null, $null, $null, $null)
    ,@(37, 35, "US Core Observation Pregnancy Intent Profile", $null, $null, $null, $null, $null, $null)
    ,@(38, 36, "US Core Observation Pregnancy Status Profile", $null, $null, $null, $null, $null, $null)
    ,@(39, 37, "US Core Observation Screening Assessment Profile", $null, $null, $null, $null, $null, $null)
    ,@(40, 38, "US Core Observation Sexual Orientation Profile", $null, $null, $null, $null, $null, $null)
    ,@(41, 39, "US Core Organization Profile", $null, $null, $null, $null, $null, $null)
    ,@(42, 40, "US Core Patient Profile", $null, $null, $null, $null, $null, $null)
    ,@(43, 41, "US Core PMO ServiceRequest Profile", $null, $null, $null, $null, $null, $null)
    ,@(44, 42, "US Core Practitioner Profile", $null, $null, $null, $null, $null, $null)
    ,@(45, 43, "US Core PractitionerRole Profile", $null, $null, $null, $null, $null, $null)
    ,@(46, 44, "US Core Procedure Profile", $null, $null, $null, $null, $null, $null)
    ,@(47, 45, "US Core Provenance Profile", $null, $null, $null, $null, $null, $null)
    ,@(48, 46, "US Core Pulse Oximetry Profile", "PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")
    ,@(49, 47, "US Core QuestionnaireResponse Profile", $null, $null, $null, $null, $null, $null)
    ,@(50, 48, "US Core RelatedPerson Profile", $null, $null, $null, $null, $null, $null)
    ,@(51, 49, "US Core Respiratory Rate Profile", "PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")
    ,@(52, 50, "US Core ServiceRequest Profile", $null, $null, $null, $null, $null, $null)
    ,@(53, 51, "US Core Simple Observation Profile", $null, $null, $null, $null, $null, $null)
    ,@(54, 52, "US Core Smoking Status Observation Profile", $null, $null, $null, $null, $null, $null)
    ,@(55, 53, "US Core Specimen Profile", $null, $null, $null, $null, $null, $null)
    ,@(56, 54, "US Core Treatment Intervention Preference Profile", $null, $null, $null, $null, $null, $null)
    ,@(57, 55, "US Core Vital Signs Profile", "PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")
)

foreach ($item in $rowsData) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    if ($item[3] -eq $null) {
        $ws.Range($ws.Cells.Item($r,3), $ws.Cells.Item($r,8)).ClearContents()
    }
    else {
        $ws.Cells.Item($r, 3).Value = $item[3]
        $ws.Cells.Item($r, 4).Value = $item[4]
        $ws.Cells.Item($r, 5).Value = $item[5]
        $ws.Cells.Item($r, 6).Value = $item[6]
        $ws.Cells.Item($r, 7).Value = $item[7]
        $ws.Cells.Item($r, 8).Value = $item[8]
    }
}

Write-Host "Update complete"